# Add two new worksheets at the end of the workbook:
#   "ODI Batting Extra"  and  "ODI Bowling Extra"
# matching the target diff (header row styled like the existing header
# rows on the other sheets; data rows as literal text / numbers).

# Helper: write a value as literal TEXT (not auto-converted to a number or
# percentage by Excel's type inference), the way entering e.g. '3632 or
# '18.21% into a cell does, then drop the quote-prefix style Excel applies
# so the cell ends up on the default (unstyled) format like the source data.
function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# Donor cell whose format (bold font, thin border, center/top alignment) is
# reused for the new header rows, so we pick up the workbook's existing
# header style instead of minting a near-duplicate one.
$headerStyleDonor = $wb.Worksheets.Item(1).Range("A1")

# ---------------------------------------------------------------------------
# Sheet 4: "ODI Batting Extra"
# ---------------------------------------------------------------------------
$wsBatExtra = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsBatExtra.Name = "ODI Batting Extra"

$wsBatExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$wsBatExtra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$wsBatExtra.Cells.Item(1, 3).Value = "NUM_4"
$wsBatExtra.Cells.Item(1, 4).Value = "NUM_6"
$wsBatExtra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$wsBatExtra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"
$headerStyleDonor.Copy()
$wsBatExtra.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

# MATCH_CODE, NUM_4, NUM_6 and PERCENT_RUNS_OF_TOTAL are stored as literal
# text in the source data (even though several look numeric/percent), so
# they are entered as text. BATTING_POSITION is a genuine number.
$batRows = @(
    @("3632", $null, $null,  $null, $null,    "NO"),
    @("3633", 2,     "5",    "2",   "18.21%", "NO"),
    @("3649", $null, $null,  $null, $null,    "NO"),
    @("3650", 2,     "12",   "1",   "46.09%", "NO"),
    @("3651", $null, $null,  $null, $null,    "NO"),
    @("3652", 2,     "0",    "1",   "5.41%",  "NO"),
    @("3705", 1,     "0",    "0",   "3.63%",  "NO"),
    @("3707", $null, $null,  $null, $null,    "NO"),
    @("3721", $null, $null,  $null, $null,    "NO"),
    @("3722", 2,     "0",    "0",   "0.56%",  "NO"),
    @("3725", 2,     "3",    "0",   "39.68%", "NO"),
    @("3730", $null, $null,  $null, $null,    "NO"),
    @("3773", 2,     "0",    "0",   "8.45%",  "NO"),
    @("3778", $null, $null,  $null, $null,    "NO"),
    @("4009", 2,     "0",    "0",   "1.36%",  "NO"),
    @("4525", $null, $null,  $null, $null,    "NO"),
    @("4528", 2,     "0",    "0",   $null,    "NO")
)

$r = 2
foreach ($row in $batRows) {
    Set-TextValue $wsBatExtra.Cells.Item($r, 1) $row[0]
    if ($null -ne $row[1]) { $wsBatExtra.Cells.Item($r, 2).Value = $row[1] }
    if ($null -ne $row[2]) { Set-TextValue $wsBatExtra.Cells.Item($r, 3) $row[2] }
    if ($null -ne $row[3]) { Set-TextValue $wsBatExtra.Cells.Item($r, 4) $row[3] }
    if ($null -ne $row[4]) { Set-TextValue $wsBatExtra.Cells.Item($r, 5) $row[4] }
    $wsBatExtra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 5: "ODI Bowling Extra"
# ---------------------------------------------------------------------------
$wsBowlExtra = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsBowlExtra.Name = "ODI Bowling Extra"

$wsBowlExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$wsBowlExtra.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$wsBowlExtra.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"
$headerStyleDonor.Copy()
$wsBowlExtra.Range("A1:C1").PasteSpecial(-4122)   # xlPasteFormats

Set-TextValue $wsBowlExtra.Cells.Item(2, 1) "3632"
Set-TextValue $wsBowlExtra.Cells.Item(2, 2) ""
Set-TextValue $wsBowlExtra.Cells.Item(2, 3) ""

Set-TextValue $wsBowlExtra.Cells.Item(3, 1) "3651"
Set-TextValue $wsBowlExtra.Cells.Item(3, 2) ""
Set-TextValue $wsBowlExtra.Cells.Item(3, 3) ""
